# Updated parser to use TokenIteratorFieldRewriterSplit.
#
# The paragraph currently reads "{m:'\t'}" held in two runs: "{m" and
# ":'\t'}". The migrated form keeps the exact same text but spreads it
# across four runs: "{", "m", ":'\t'", "}". Word's object model has no
# direct "split this run" verb (and it merges adjacent runs that end up
# with identical formatting when the document is saved), so the split
# points are forced by dropping a zero-length bookmark at each boundary
# and immediately deleting it again: adding the bookmark forces the run
# to break at that character offset, and removing the bookmark again
# leaves the run break behind without adding any visible formatting.

$d = $word.ActiveDocument

# Locate the paragraph holding the "{m:'\t'}" field text (literal
# backslash-t, not a tab character). Compare with StartsWith rather than
# equality since Range.Text includes the trailing paragraph mark.
$needle = "{m:'\t'}"
$para = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith($needle)) {
        $para = $p
    }
}

$start = $para.Range.Start

# Character offsets (relative to the paragraph start) at which a run
# break is required so the text splits into "{" / "m" / ":'\t'" / "}":
#   offset 1 -> between "{" and "m"
#   offset 2 -> between "m" and ":'\t'"
#   offset 7 -> between ":'\t'" and "}"
$splitOffsets = 1, 2, 7

$i = 0
foreach ($off in $splitOffsets) {
    $i = $i + 1
    $pos = $start + $off
    $bmk = $d.Range($pos, $pos)
    $name = "m2docSplit" + $i
    $d.Bookmarks.Add($name, $bmk)
}

$i = 0
foreach ($off in $splitOffsets) {
    $i = $i + 1
    $name = "m2docSplit" + $i
    $d.Bookmarks($name).Delete()
}
